$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / Column A: tab label changed from "CasesTab" to "ProgramsTab"
$ws.Range("A2").Value = "ProgramsTab"

# Row 2 / Column B: same Cypher query text, with a one-space whitespace fix
# on the "PubMed ID" coalesce line (matches the authoritative diff).
$ws.Range("B2").Value = @"
MATCH (ss:study_subject)
Match (s:study)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, 
collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, 
collect(DISTINCT f) AS files
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
OPTIONAL MATCH (ss)<-[:program_of_institution]-(p)
OPTIONAL MATCH (p)<-[:of_arm]-(a)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
RETURN DISTINCT
       coalesce (p.program_acronym, '')as ``Program Code``,
       coalesce( p.program_id , '')as ``Program ID``,
       coalesce (p.program_name , '' )as ``Program Name``,
       coalesce(p.start_date, '') as ``Start Date``,
       coalesce (p.end_date, '') as ``End Date``,
      coalesce(p.pubmed_id, '') as ``PubMed ID``,
       count(distinct s) As ``Number of Arms``,
       count(distinct ss) as ``Associated Cases``
       order By ``Program Code``
"@

# Selection moved from B2 to C5
$ws.Range("C5").Select()
